# Add a 2nd test worksheet's worth of data to the existing (empty) sheet:
#  - rename the sheet from "Лист1" to "AddCustomerTest"
#  - fill a 3x2 block of header + data (firstname/lastname/postcode, Katya/Smith/ab214c)
#  - leave the selection on C2, matching the saved Excel session state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "AddCustomerTest"

$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

$ws.Range("A2").Value = "Katya"
$ws.Range("B2").Value = "Smith"
$ws.Range("C2").Value = "ab214c"

# Match the saved selection (active cell C2) from the target workbook.
$ws.Range("C2").Select() | Out-Null
